$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row (A1:D1); columns E1:I1 are no longer used
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"
$ws.Range("E1:I1").Clear()

# Update the single remaining data row (row 2); keep only columns A:D
$ws.Range("A2").Value = 42
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 3144
$ws.Range("D2").Value = 441.2654347419739
$ws.Range("E2:I2").Clear()

# Remove rows 3 and 4, which are no longer needed
$ws.Rows("3:4").Delete()
